$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-67) holds a date serial value that was bumped from
# 45203 (2023-10-02) to 45205 (2023-10-04) for every row in the sheet.
$ws.Range("C2:C67").Value = 45205
